$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the small histogram bucket table in columns H (bucket) and I (count),
# extending rows 67-77.
$ws.Range("H67").Value = -5
$ws.Range("I67").Formula = '=COUNTIFS($H$2:$H$65,H67)'

$row = 68
for ($bucket = -4; $bucket -le 5; $bucket++) {
    $ws.Cells.Item($row, 8).Value = $bucket
    $ws.Cells.Item($row, 9).Formula = '=COUNTIFS($H$2:$H$65,H' + $row + ')'
    $row++
}

# Update the active selection on the sheet (canvas small fix)
$ws.Range("A21:H21").Select()
